$d = $word.ActiveDocument
$wns = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

function Replace-RunText($find, $replace, $rPrXml) {
    # Locate the plain-text range for $find and replace its contents with a
    # freshly built <w:r> (wrapped in a throw-away <w:p> so InsertXML parses
    # it correctly) carrying the same run formatting. This performs a
    # surgical text-only substitution without disturbing any neighbouring
    # (e.g. empty) runs, matching how a translation/CAT tool would normally
    # edit the underlying OOXML.
    $text = $d.Content.Text
    $idx = $text.IndexOf($find)
    if ($idx -lt 0) {
        return $false
    }
    $rng = $d.Range($idx, $idx + $find.Length)
    $escaped = $replace.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $spacePreserve = ""
    if ($replace -ne $replace.Trim()) {
        $spacePreserve = " xml:space=`"preserve`""
    }
    $xml = "<w:p xmlns:w=`"$wns`"><w:r>$rPrXml<w:t$spacePreserve>$escaped</w:t></w:r></w:p>"
    $rng.InsertXML($xml) | Out-Null
    return $true
}

# 1. Main page title (Heading1) - no neighbouring empty run, so a plain
#    Find/Replace of just the first occurrence is safe and simplest.
$d.Content.Find.Execute("Play Chilli Heat Free: Review of Pragmatic Play Slot", $true, $false, $false, $false, $false, $true, 1, $false, "Play Chilli Heat Free: Exciting Online Slot Game", 1)

# 2-5. "What we like" bullet list items
Replace-RunText "Engaging gameplay with straightforward controls" "Engaging Mexican cuisine theme" "" | Out-Null
Replace-RunText "Cartoonish graphics and fiery red background" "Two bonus games for extra winning potential" "" | Out-Null
Replace-RunText "Two bonus games for significant payouts" "Wide range of betting options" "" | Out-Null
Replace-RunText "Wide range of betting options from €0.25 to €125 per spin" "Decent RTP of 96.5%" "" | Out-Null

# 6. "What we don't like" bullet list item
Replace-RunText "Only the Free Spins feature offers richer reels" "Limited number of free spins" "" | Out-Null

# 7. Bolded title repeated near the end of the document
Replace-RunText "Play Chilli Heat Free: Review of Pragmatic Play Slot" "Play Chilli Heat Free: Exciting Online Slot Game" "<w:rPr><w:b/></w:rPr>" | Out-Null

# 8. Italicised meta description
Replace-RunText "Chilli Heat is a fun and engaging slot game based on Mexican cuisine featuring two bonus games for significant payouts and a 96.5% RTP. Play for free now!" "Read our review of Chilli Heat and play this exciting online slot game for free. Experience Mexican cuisine and win big!" "<w:rPr><w:i/></w:rPr>" | Out-Null
